$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Save" header in H1, reusing the same formatting as the
# existing header row (bold font, border, centered) by copying G1's
# format rather than rebuilding it (keeps the shared style index stable).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the corresponding numeric value for H2 (no special style, like B2:G2)
$ws.Range("H2").Value = 1
